# Adds a new weekly block of 3 price rows (Especial/Primera/Segunda) for
# "Mango" at "Terminal La Palmera de La Serena" with a new date/origin,
# inserted immediately above the existing 2021-03-04 (serial 44263) block.
# Inserting the rows pushes the 91 existing data rows (853:943) down to
# (856:946), which is why the sheet's used range grows from T943 to T946.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above row 853 (shifts 853:943 -> 856:946).
$ws.Rows("853:855").Insert()

$qualities = @("Especial", "Primera", "Segunda")

for ($i = 0; $i -lt 3; $i++) {
    $r = 853 + $i

    $ws.Cells.Item($r, 1).Value = 8
    $ws.Cells.Item($r, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 4).Value = 44858
    $ws.Cells.Item($r, 5).Value = 4
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100108
    $ws.Cells.Item($r, 8).Value = "Tropicales y subtropicales"
    $ws.Cells.Item($r, 9).Value = 100108002
    $ws.Cells.Item($r, 10).Value = "Mango"
    $ws.Cells.Item($r, 11).Value = "Sin especificar"
    $ws.Cells.Item($r, 12).Value = $qualities[$i]
    $ws.Cells.Item($r, 13).Value = 512
    $ws.Cells.Item($r, 14).Value = 7500
    $ws.Cells.Item($r, 15).Value = 8000
    $ws.Cells.Item($r, 16).Value = 7750
    $ws.Cells.Item($r, 17).Value = "`$/bandeja 4 kilos"
    $ws.Cells.Item($r, 18).Value = "Brasil"
    $ws.Cells.Item($r, 19).Value = 1938
    $ws.Cells.Item($r, 20).Value = 4
}
